# Apply BOM changes: added R43 (4.7k), R44 (300k), R41/R42 moved to 10k group,
# and replaced the dual op-amp (U1, TL972IDR, SOIC-8) with a quad op-amp
# (TL974IDR, SOIC-14) to restore mono operation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 (ID 9): 2.20k resistors - R12 and R19 removed, leaving R1, R15
$ws.Range("B10").Value = "R1, R15"
$ws.Range("D10").Value = 2

# Row 11 (ID 10): 4.7k resistor - R43 added
$ws.Range("B11").Value = "R10, R43"
$ws.Range("D11").Value = 2

# Row 12 (ID 11): 750 resistors - R13 and R18 removed, leaving R11, R14, R17
$ws.Range("B12").Value = "R11, R14, R17"
$ws.Range("D12").Value = 3

# Row 16 (ID 15): 10k resistors - R41 and R42 added
$ws.Range("B16").Value = "R26, R27, R30, R33, R41, R42"
$ws.Range("D16").Value = 6

# Row 20 (ID 19): 300k resistors - R44 added
$ws.Range("B20").Value = "R36, R37, R44"
$ws.Range("D20").Value = 3

# Row 22 (ID 21): 1k resistors - R41 and R42 removed, leaving R6, R7, R8, R9
$ws.Range("B22").Value = "R6, R7, R8, R9"
$ws.Range("D22").Value = 4

# Row 23 (ID 22): U1 swapped from dual op-amp (TL972IDR, SOIC-8) to
# quad op-amp (TL974IDR, SOIC-14) so the circuit can be restored to mono.
$ws.Range("C23").Value = "SOIC-14"
$ws.Range("E23").Value = "TL974IDR"
$ws.Range("F23").Value = "TL974IDR"

# Restore selection/scroll to match the saved view state (A1 top-left, D11 active).
$ws.Range("A1").Select()
$ws.Range("D11").Select()
